$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update nutrient labels (column B) to restore correct order ---
$ws.Range("B7").Value = "COLEST"
$ws.Range("B8").Value = "CALCIO"
$ws.Range("B9").Value = "AGTRANS"
$ws.Range("B10").Value = "AGSAT"
$ws.Range("B11").Value = "AGPOLI"
$ws.Range("B18").Value = "PIRIDOXAMINA"
$ws.Range("B19").Value = "NIACINA"

# --- Update Final Value (D) and Target Value (E) columns with corrected heuristic results ---
# (Initial Value column C also reassigned for rows where row order changed)
$ws.Range("C2").Value = 1404.55
$ws.Range("D2").Value = 4332.95
$ws.Range("E2").Value = 2989.87
$ws.Range("C3").Value = 194.79
$ws.Range("D3").Value = 805.3
$ws.Range("E3").Value = 1644.4285
$ws.Range("C4").Value = 90.52
$ws.Range("D4").Value = 111.79
$ws.Range("E4").Value = 298.987
$ws.Range("C5").Value = 35.89
$ws.Range("D5").Value = 81.74
$ws.Range("E5").Value = 448.4805000000001
$ws.Range("C6").Value = 29.47
$ws.Range("D6").Value = 34.27
$ws.Range("E6").Value = 31
$ws.Range("C7").Value = 122.22
$ws.Range("D7").Value = 261.75
$ws.Range("E7").Value = 300
$ws.Range("C8").Value = 219.6
$ws.Range("D8").Value = 720
$ws.Range("E8").Value = 868
$ws.Range("C9").Value = 0.9
$ws.Range("D9").Value = 1.25
$ws.Range("E9").Value = 29.89870000000001
$ws.Range("C10").Value = 11.01
$ws.Range("D10").Value = 29.47
$ws.Range("E10").Value = 298.987
$ws.Range("C11").Value = 11.47
$ws.Range("D11").Value = 20.61
$ws.Range("E11").Value = 179.3922
$ws.Range("C12").Value = 2623.22
$ws.Range("D12").Value = 3059.56
$ws.Range("E12").Value = 1
$ws.Range("C13").Value = 2241.43
$ws.Range("D13").Value = 3721.76
$ws.Range("E13").Value = 3510
$ws.Range("C14").Value = 13.55
$ws.Range("D14").Value = 17.37
$ws.Range("E14").Value = 6.8
$ws.Range("C15").Value = 309.36
$ws.Range("D15").Value = 326.25
$ws.Range("E15").Value = 303
$ws.Range("C16").Value = 0.7
$ws.Range("D16").Value = 1.2
$ws.Range("E16").Value = 0.9
$ws.Range("C17").Value = 0.95
$ws.Range("D17").Value = 2.15
$ws.Range("E17").Value = 1
$ws.Range("C18").Value = 0.52
$ws.Range("D18").Value = 0.96
$ws.Range("E18").Value = 1.1
$ws.Range("C19").Value = 11.74
$ws.Range("D19").Value = 16.52
$ws.Range("E19").Value = 11.5
$ws.Range("C20").Value = 4.06
$ws.Range("D20").Value = 5.76
$ws.Range("E20").Value = 2
$ws.Range("C21").Value = 3.74
$ws.Range("D21").Value = 51.94
$ws.Range("E21").Value = 66.1
$ws.Range("C22").Value = 41.23
$ws.Range("D22").Value = 661.71
$ws.Range("E22").Value = 560
$ws.Range("C23").Value = 1.39
$ws.Range("D23").Value = 1.56
$ws.Range("E23").Value = 0.7
$ws.Range("C24").Value = 481.27
$ws.Range("D24").Value = 552.47
$ws.Range("E24").Value = 322
$ws.Range("C25").Value = 989.78
$ws.Range("D25").Value = 1406.34
$ws.Range("E25").Value = 649
$ws.Range("C26").Value = 16.71
$ws.Range("D26").Value = 19.54
$ws.Range("E26").Value = 8
